$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Port B Interrupt now fully implemented: IST (actual) points for the
# "Interruptfunktion" row go from 3 to 5, and the old comment explaining the
# partial state is no longer needed.
$ws.Range("D15").Value = 5
$ws.Range("F15").ClearContents()

# Move the active selection to reflect where the author was last working.
$ws.Range("G28").Select()

$wb.Save()
